$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: 16.42578125 -> 15.7109375 (closest achievable via ColumnWidth granularity)
$ws.Columns.Item(1).ColumnWidth = 14.83

# Updated simulation values in column A (rod contact with ground example)
$ws.Range("A7").Value = 0.040000393448394019
$ws.Range("A9").Value = 0.000080904185094077174
$ws.Range("A10").Value = 0.060000402724620663
$ws.Range("A12").Value = 0.00022893671078018558
$ws.Range("A13").Value = 0.079999919458515353
$ws.Range("A15").Value = 0.0004329364119163195
$ws.Range("A16").Value = 0.099998902518024999
$ws.Range("A18").Value = 0.00068385998607482616
$ws.Range("A19").Value = 0.11999734958525113
$ws.Range("A21").Value = 0.00097437919335646695
$ws.Range("A22").Value = 0.13999527941231149
$ws.Range("A24").Value = 0.001298554321017023
$ws.Range("A25").Value = 0.15999272093643743
$ws.Range("A27").Value = 0.0016515694584630218
$ws.Range("A28").Value = 0.17998970677200701
$ws.Range("A30").Value = 0.0020295178313289712
$ws.Range("A31").Value = 0.19998626947172521
$ws.Range("A33").Value = 0.0024292276486627255
$ws.Range("A34").Value = 0.21998243951985058
$ws.Range("A36").Value = 0.0028481207103325779
$ws.Range("A37").Value = 0.23997824439218543
$ws.Range("A39").Value = 0.0032840974759685551
$ws.Range("A40").Value = 0.2599737082593247
$ws.Range("A42").Value = 0.0037354434734995414
$ws.Range("A43").Value = 0.27996885206651662
$ws.Range("A45").Value = 0.0042007528753908296
$ws.Range("A46").Value = 0.29996369382498672
$ws.Range("A48").Value = 0.0046788658355207896
$ws.Range("A49").Value = 0.31995824901507114
$ws.Range("A51").Value = 0.0051688167929343463
$ws.Range("A52").Value = 0.33995253104366951
$ws.Range("A54").Value = 0.0056697914376085306
$ws.Range("A55").Value = 0.35994655172565365
$ws.Range("A57").Value = 0.0061810904195572927
$ws.Range("A58").Value = 0.37994032177636933
$ws.Range("A60").Value = 0.006702098183261015
$ws.Range("A61").Value = 0.39993385131379067
$ws.Range("A63").Value = 0.0072322555379183652
$ws.Range("A64").Value = 0.41992715037654027
$ws.Range("A66").Value = 0.0077710347405684057
$ws.Range("A67").Value = 0.43992022946935172
$ws.Range("A69").Value = 0.0083179159811449608
$ws.Range("A70").Value = 0.45991310015158693
$ws.Range("A72").Value = 0.0088723642209960398
$ws.Range("A73").Value = 0.479905775687682
$ws.Range("A75").Value = 0.0094338053521087137
$ws.Range("A76").Value = 0.49989827178114837
$ws.Range("A78").Value = 0.010001600613914157
$ws.Range("A79").Value = 0.51989060741600757
$ws.Range("A81").Value = 0.010575018126738845
$ws.Range("A82").Value = 0.53988280583099546
$ws.Range("A84").Value = 0.011153200272224995
$ws.Range("A85").Value = 0.55987489565185555
$ws.Range("A87").Value = 0.011735125465595434
$ws.Range("A88").Value = 0.57986691220432629
$ws.Range("A90").Value = 0.012319562614175328
$ws.Range("A91").Value = 0.5998588990229573
$ws.Range("A93").Value = 0.012905016229881178
$ws.Range("A94").Value = 0.61985090955529687
$ws.Range("A96").Value = 0.013489659745802849
$ws.Range("A97").Value = 0.63984300903196145
$ws.Range("A99").Value = 0.014071254059820612
$ws.Range("A100").Value = 0.65983527642218476
$ws.Range("A102").Value = 0.014647047667795883
$ws.Range("A103").Value = 0.67982780630849993
$ws.Range("A105").Value = 0.015213653925700402
$ws.Range("A106").Value = 0.69982071037271587
$ws.Range("A108").Value = 0.015766899957368023
$ws.Range("A109").Value = 0.7198141179565607
$ws.Range("A111").Value = 0.01630164045696484
$ws.Range("A112").Value = 0.73980817479519001
$ws.Range("A114").Value = 0.01681152806702611
$ws.Range("A115").Value = 0.75980303844460018
$ws.Range("A117").Value = 0.01728873007613246
$ws.Range("A118").Value = 0.77979886801858467
$ws.Range("A120").Value = 0.017723578793439481
$ws.Range("A121").Value = 0.79979580443883058
$ws.Range("A123").Value = 0.01810414002438342
$ws.Range("A124").Value = 0.81979393520974508
$ws.Range("A126").Value = 0.018415680484108535
$ws.Range("A127").Value = 0.83979323433963082
$ws.Range("A129").Value = 0.018640010627159458
$ws.Range("A130").Value = 0.85979346280382929
$ws.Range("A132").Value = 0.018754674138618792
$ws.Range("A133").Value = 0.87979400691193788
$ws.Range("A135").Value = 0.018731949165536704
$ws.Range("A136").Value = 0.8997936196259666
$ws.Range("A138").Value = 0.01853761933592199
$ws.Range("A139").Value = 0.91979001105365155
$ws.Range("A141").Value = 0.018129465064807206
$ws.Range("A142").Value = 0.93977920567652906
$ws.Range("A144").Value = 0.017455418518494869
$ws.Range("A145").Value = 0.95975454043038821
$ws.Range("A147").Value = 0.016451321007359018
$ws.Range("A148").Value = 0.97970511237447111
$ws.Range("A150").Value = 0.01503822389546162
$ws.Range("A151").Value = 0.99961338731416638
$ws.Range("A153").Value = 0.013119192093683009
